$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Modify existing shared-string cells ---
$ws.Range("D9").Value = '<p><strong>Note</strong>: <em>This quest might take a moment, longer then two hours, because you need the gold dust which you get from enchanting and then disenchanting.</em></p><p>We know how to do quests now, so lets continue doing quests to get an item that will help you – especially if you are crafting items with Accuracy enchantments attached – level Accuracy rather quickly.</p><p>This right here shows you how two features of the game can work together to progress your character: Quest items that you get from doing quests that increase some aspect about your character. Enchanting, Crafting, Quests and Skills</p><p>As we progress we will open more systems that also work together with existing systems you have learned to continue to progress your character.</p>'
$ws.Range("D10").Value = '<p><strong>Note</strong>: <em>This quest can take two plus hours to complete. I recommend exploration while crafting and enchanting.</em></p><p>It’s time we get off the <a href="/information/planes" target="_blank">plane</a>: Surface. To do this we just have to kill a single creature: Labyrinth Fiend until the Key of Labyrinth drops.</p><p>This will allow us to use a new feature called: <a href="/information/traverse" target="_blank">Traverse</a>.</p><p>Traverse allows us to move to another plane of existence where we can fight monsters and interact with new features. For example, when we eventually get to Hell, those <a href="/information/random-enchants" target="_blank">Green Uniques</a> you have been getting – there is an NPC that allows you to buy better ones and even re-roll them and move the enchantments to other items. How fun!</p><p>Some Planes make you weaker and the creatures much, much stronger and also give you bonuses for killing creatures of that plane. As you go past Dungeons plane you will find your gear needs to be better with a set of enchantments that play to your classes strengths as well as the levels.</p>'
$ws.Range("O12").Value = 'child-of-shards'

# --- Row 12 P/Q value changes ---
$ws.Range("P12").Value = 2
$ws.Range("Q12").Value = 2

# --- Column E width ---
$ws.Columns.Item(5).ColumnWidth = 1851.15


# --- New row 13 ---
$ws.Range("A13").Value = 26
$ws.Range("B13").Value = 'Farmers in the fields'
$ws.Range("C13").Value = 'You spent your coin, weather you won or lost is up for debate, but you were able to get your mind off the horrors of the INN.<br /> <br /> You wonder for days passing through farm fields, stopping to speak to a few farmers who tell you of their lord, their lady, their king, their queen.<br /> <br /> Seems like there is a lot of a kingdoms around with various lords and ladies and since you haven’t seen The Guide since he last appear coughing up blood and looking half dead, you decide to investigate this kingdom concept.<br /> <br /> Can you just settle any where?<br /> <br /> “Yup, that’s how it works.” Stats a farmer before he walks off back to his fields, leaving you to stand alone in the warm fall day.<br /> <br /> Winter is coming after all, the nights are getting colder and with every one murdered by something, a mysterious shadow perhaps, you need a new home.'
$ws.Range("D13").Value = '<p>Lets settle a <a href="/information/kingdoms" target="_blank">kingdom</a>. Kingdoms become important later on in the game and require the player to be vigilant because after 7 days from settling your first kingdom, it is open to everyone to attack and take.</p><p>One of the systems that kingdoms impact is a later system called <a href="/information/gems" target="_blank">Atonement</a>, which includes adding sockets and gems to your gear to build up resistance and elemental damage for end game <a href="/information/raids" target="_blank">Raids</a> and associated critters as well as <a href="/information/player-vs-player" target="_blank">PVP</a>.</p><p>To socket gems on to items you will need something called Gold Bars, which you can only get by owning one or more kingdoms.</p><p>For now, lets focus on settling a kingdom and building up a bit of defense, you can do this with gold or with kingdom resources which slowly accumulate once per hour.</p><p><strong>Caution</strong>: There is an option in your kingdom to allow you to buy population. Your kingdom, based on the farm level, has a maximum population. You can buy population beyond this. But, should the hour be up and you have not brought that population down to your current max, an NPC: The Old Man, will stomp around in chat. He will first take payment from your Gold Bars. Since you won’t have any he next takes what he can from your treasury, if there''s not enough or no gold he will take the rest from your personal gold. If you have no gold or not enough he will demolish your kingdom.</p><p><a href="/information/npc-kingdoms" target="_blank"><strong>Yellow Kingdoms</strong></a>: You might see Yellow Kingdoms on the map, these are abandoned or neglected kingdoms. Kingdoms will loose morale if you do not visit them enough, by actually going to them once a month at least. This is called walking the kingdoms. Failure to do this once per 30 days can see your kingdom become abandoned, these can be bought by other players from The Old Man.</p><p><strong>Red Kingdoms:</strong> These are enemy kingdoms. Simple as that. You can attack them and later on, when you have settlers, take kingdoms. Some Alchemy items can be dropped on them to do wild amounts of damage, making the process of taking them with settlers much easier.</p>'
$ws.Range("E13").Value = '<p>- Find a spot on the map you can move to, could be Surface or Labyrinth since you also have access to Labyrinth. The only plane you cannot settle on is Purgatory.</p><p>- This location must be 1 space away from any location or port – there must be one space between your kingdom and a port for example, in all directions.</p><p>- Under the map click Settle Kingdom.</p><p>- Name your kingdom.</p><p>- Click settle.</p><p>- Boom done!</p><p>You could also look to see if there is a yellow kingdom around. If so click it and purchase it, it will be free. You can rename it from the kingdom section.</p><p>Your very first, even if you loose all your kingdoms, kingdom is always fee, every kingdom after that, including settling, is 10,000 Gold x Kingdom amount. Ie, 6 kingdoms: 60K god your 7<sup>th</sup> kingdom.</p><p>- Click the kingdoms tab.</p><p>- Click your kingdom name.</p><p>Here you will see all the info about this kingdom, including buildings and their levels as well as other actions.</p><p>You can click on buildings and the unit names to then upgrade or, if they have been damaged in someone attacking, repair them. You can do this with gold or resources.</p><p>The requirements wants a sum total of the levels of the your kingdom – this means total level of all your buildings. You could do one or multiple, I suggest the farm, walls and barracks but I’ll let you explore. Buildings and unit can take minutes or hours or even days to process.</p><p>The other requirement for units, is the same concept – the sum total of one or more sets of units.</p>'
$ws.Range("F13").Value = '<p>- Select Map Movement from the actions section</p><p>- Find a spot on the map you can move to, could be Surface or Labyrinth since you also have access to Labyrinth. The only plane you cannot settle on is Purgatory.</p><p>- This location must be 1 space away from any location or port – there must be one space between your kingdom and a port for example, in all directions.</p><p>- Under the map tap Settle Kingdom.</p><p>- Name your kingdom.</p><p>- tap settle.</p><p>- Boom done!</p><p>You could also look to see if there is a yellow kingdom around. If so teleport to it by selecting a NPC Kingdom from the NPC Kingdoms Drop down. Purchase it. It will be free. You can rename it from the kingdom section.</p><p>Your very first, even if you loose all your kingdoms, kingdom is always fee, every kingdom after that, including settling, is 10,000 Gold x Kingdom amount. Ie, 6 kingdoms: 60K god your 7<sup>th</sup> kingdom.</p><p>- Tap the kingdoms tab.</p><p>- Tap your kingdom name.</p><p>Here you will see all the info about this kingdom, including buildings and their levels as well as other actions.</p><p>You can tap on buildings and the unit names to then upgrade or, if they have been damaged in someone attacking, repair them. You can do this with gold or resources.</p><p>The requirements wants a sum total of the levels of the your kingdom – this means total level of all your buildings. You could do one or multiple, I suggest the farm, walls and barracks but I’ll let you explore. Buildings and unit can take minutes or hours or even days to process.</p><p>The other requirement for units, is the same concept – the sum total of one or more sets of units.</p>'
$ws.Range("AC13").Value = 1
$ws.Range("AD13").Value = 30
$ws.Range("AE13").Value = 100
$ws.Range("AR13").Value = 1000000
$ws.Range("AS13").Value = 100

# --- New row 14 ---
$ws.Range("A14").Value = 27
$ws.Range("B14").Value = 'Blasting shadows'
$ws.Range("C14").Value = 'It has been months, the winter has come and gone and the fields are starting to grow bountiful food. You start the plans for the next leg of your building adventure.<br /> <br /> The back of your mind is still haunted by the visions of death, most of them have faded into the distance, into the shadows.<br /> <br /> You look up from your papers and plans and see a shadow cross the dimly lit room. You stand and ready your magics, you illuminate the room and the light blasts through the room and out the windows into the darkness of the night.<br /> <br /> The people below look up and see the room you are in brightly lit up, such that they must shield their eyes.<br /> <br /> The shadow is gone, but the fear still stands in the back of your mind.'
$ws.Range("D14").Value = '<p><strong>Note</strong>: <em>This quest will take two or more hours to complete. Utilize exploration to level your class ranks and mercenaries while you craft and then enchant and then – if it’s not better then what you have equiped – disenchant it.</em></p><p>Lets learn about <a href="/information/kingdom-passive-skills" target="_blank">Kingdom Passives</a>, these are a special set of skills that take real world hours to train a single level which increases as the level of the passive increases.</p><p>Each child skill will require it’s parent to be at a specific level to unlock the ability to level it, for this reason it’s advised players go after the Passive they want, in our case – we are going to head for Goblin Coin Bank and while we do that were going to work on our crafting and enchanting.</p><p>Kingdom Passives help you gain bonuses towards all kingdoms as well as unlock new buildings and units. For example the thing we want: Goblin Coin Bank, will allow you to construct and level this building to take 2 billion gold and store it in your kingdom as 1 gold bar. A single kingdom have 1000 of these bars. Gold bars are then used in things like Sockets and the <a href="/information/goblin-shop" target="_blank">Goblin Shop</a>.</p>'
$ws.Range("E14").Value = '<p>- Click Character Sheet tab</p><p>- Under skills section is a tab called Kingdom Passives</p><p>- Click the top Skill and begin training it to level 1.</p><p>These passives will level while you are logged out.</p>'
$ws.Range("F14").Value = '<p>- Tap Character Sheet tab</p><p>- Select Skill management from the drop down</p><p>- Under skills section is a tab called Kingdom Passives</p><p>- Tap the top Skill and begin training it to level 1.</p><p>These passives will level while you are logged out.</p>'
$ws.Range("H14").Value = 'Enchanting'
$ws.Range("I14").Value = 25
$ws.Range("J14").Value = 'Ring Crafting'
$ws.Range("K14").Value = 5
$ws.Range("AD14").Value = 60
$ws.Range("AE14").Value = 500
$ws.Range("AF14").Value = 'Kingdom Management'
$ws.Range("AG14").Value = 1
$ws.Range("AR14").Value = 2000000
$ws.Range("AS14").Value = 100

# --- New row 15 ---
$ws.Range("A15").Value = 28
$ws.Range("B15").Value = 'Sliding into despair'
$ws.Range("C15").Value = 'You have heard reports of the wars that have been raging between the various kingdoms in your area. You are worried that your kingdom could be next and you hurry to build the defenses around your kingdom fearing that it’s only a matter of time, slowly drifting into the shadows of delusion while the shadow haunts you.'
$ws.Range("D15").Value = '<p>We are going to continue training out kingdom passives. This might take a couple days. The reason is, I want you to train the skill: Building Management to level 2. To do that you need to train the parent skill to level 4. This alone is another 4 hours, on top of 4 hours to get the required skill to level 2.</p><p>In the mean time lets focus on crafting and enchanting.</p>'
$ws.Range("E15").Value = '<p>Repeat the steps of the previous quest. The required passive will become unlocked when your parent skill is level 4, so come back often to check the status of the passive.</p>'
$ws.Range("F15").Value = '<p>Repeat the steps of the previous quest. The required passive will become unlocked when your parent skill is level 4, so come back often to check the status of the passive.</p>'
$ws.Range("AD15").Value = 75
$ws.Range("AE15").Value = 2000
$ws.Range("AF15").Value = 'Building Management'
$ws.Range("AG15").Value = 2
$ws.Range("AR15").Value = 2000000
$ws.Range("AS15").Value = 100

# --- New row 16 ---
$ws.Range("A16").Value = 29
$ws.Range("B16").Value = 'Goblins and Armies'
$ws.Range("C16").Value = 'Kingdoms have fallen, the death of the lords around you, the blood of their soldiers soak the lands. You are sitting in your throne room, in the darkness, in the shadows of your mind. Haunted by death, haunted by the people you rule, the people you lead, the people you protect.<br /> <br /> The days have gone by, the plans have the been seeded and they flower into your ambitious outcomes.<br /> <br /> You have heard of the goblins who run special banks that allow the super rich to store their vast wealth and exchange them for goods and services at special shops. Now might be the time to distract our self from the death and carnage and look into this concept.<br /> <br /> Wonder if The Helpless Goblin might have anything of use.'
$ws.Range("D16").Value = '<p><strong>Note:</strong> <em>This quest will take up to 4 or more hours to complete</em></p><p>Now lets go and get out Goblin Coin Bank to level 1 and this will unlock the new Goblin Bank Building that you can upgrade to level 5.</p><p>Once the building is level 5 anew action button will appear where you can manage the gold bars in your kingdom.</p><p>One thing to note: Kingdom Defense</p><p>Your kingdom defense will protect you, along with your units, against incoming attacks as well as alchemical items that get dropped on your kingdom.</p><p>This defense is made up of the level of your wall, how much treasury you have, how many gold bars you have and other factors. The more gold in the treasury and the more bars in the bank, the higher the defense, the wall will also help with defense as you level it.</p>'
$ws.Range("E16").Value = '<p>At this point we are just repeating the steps of the past and leveling some other aspects of our character.</p>'
$ws.Range("F16").Value = '<p>At this point we are just repeating the steps of the past and leveling some other aspects of our character.</p>'
$ws.Range("H16").Value = 'Armour Crafting'
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 'Enchanting'
$ws.Range("K16").Value = 50
$ws.Range("R16").Value = 2
$ws.Range("S16").Value = 10
$ws.Range("W16").Value = 'Leading Armies'
$ws.Range("AD16").Value = 100
$ws.Range("AE16").Value = 3000
$ws.Range("AF16").Value = 'Goblin Coin Bank'
$ws.Range("AG16").Value = 1
$ws.Range("AR16").Value = 5000000
$ws.Range("AS16").Value = 100

